# Adds two new slides at the end of the deck, each containing a single
# "Content Placeholder" text box (no title) with a link related to the
# exercise repository, per the "Added links in slides." commit.

$p = $ppt.ActivePresentation

# EMU -> point conversion factor used by the Shape position/size properties.
$emuPerPt = 12700

# --- Slide 4: plain clone-able repo URL -----------------------------------
$s4 = $p.Slides.Add(4, [PpSlideLayout]::ppLayoutText)

# The layout comes with a Title placeholder we don't want on this slide.
$s4.Shapes.Item(1).Delete()

$ph4 = $s4.Shapes.Item(1)
$ph4.Left = 457200 / $emuPerPt
$ph4.Top = 2514600 / $emuPerPt
$ph4.Width = 8229600 / $emuPerPt
$ph4.Height = 1600200 / $emuPerPt

$tr4 = $ph4.TextFrame.TextRange
$tr4.Text = "https://github.com/mattphotonman/2014-03-17-nyu-exercises"
$tr4.ParagraphFormat.Bullet.Type = [PpBulletType]::ppBulletNone

# --- Slide 5: "git clone ..." command line ---------------------------------
$s5 = $p.Slides.Add(5, [PpSlideLayout]::ppLayoutText)
$s5.Shapes.Item(1).Delete()

$ph5 = $s5.Shapes.Item(1)
$ph5.Left = 457200 / $emuPerPt
$ph5.Top = 2514600 / $emuPerPt
$ph5.Width = 8229600 / $emuPerPt
$ph5.Height = 1447800 / $emuPerPt

$tr5 = $ph5.TextFrame.TextRange
$tr5.Text = "g"
$tr5.InsertAfter("it")
$tr5.InsertAfter(" clone https")
$tr5.InsertAfter("://github.com")
$tr5.InsertAfter("/<username>/2014-03-17-nyu-exercises.git")

$tr5b = $ph5.TextFrame.TextRange
$tr5b.ParagraphFormat.Bullet.Type = [PpBulletType]::ppBulletNone

$p.Save()
